$wb = $excel.ActiveWorkbook

# --- Sheet "Version & History": insert the V1.5 changelog row ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows(10).Insert()
$ws1.Range("A10:E10").Style = $ws1.Range("A9:E9").Style

$ws1.Range("A10").Value = "V1.5"
$ws1.Range("B10").Value = "Fixed both accelerator and brake pedal resolutions and also fixed the steering wheel angle resolution."
$ws1.Range("C10").Value = "Bertalan " + [char]0x00C1 + "d" + [char]0x00E1 + "m"
$ws1.Range("D10").Value = 42803
$ws1.Range("E10").Value = "Draft version"

# --- Sheet "CommunicationMatrix": fix pedal/steering resolutions & drop the Turn Signals row ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("E8").Value = 1
$ws2.Range("E9").Value = 1
$ws2.Range("E10").Value = 1

# Remove the (temporary/future feature) Headlight row; Turn Signals row above effectively
# becomes the last "Output" style row once the strings collapse together.
$ws2.Rows(15).Delete()

$ws2.Range("B14").Value = "Headlight"
$ws2.Range("G14").Value = "ON: TRUE" + [char]10 + "OFF: FALSE"
$ws2.Range("K14").Value = "This signal will set the visual indicator for the headlight state"
$ws2.Rows(14).RowHeight = 23.25
